# The "Metadata" sheet lists FHIR ImplementationGuide metadata as
# Property/Value pairs in columns A/B. This deploy updates the published
# Status from "active" to "draft" and bumps the Date to the new
# publication timestamp.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

$statusLabel = $ws.Range("A1:A20").Find("Status")
$statusLabel.Offset(0, 1).Value = "draft"

$dateLabel = $ws.Range("A1:A20").Find("Date")
$dateLabel.Offset(0, 1).Value = "2023-08-01T16:12:28+00:00"
